$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.37%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.25%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.482"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08026"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.40%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.016"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.31%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9537"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.78%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.560"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.40%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1144"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.41%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1895"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.43%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "10.69"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "26.59%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09921"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.49%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04818"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12.26%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1065"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.23%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001267"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.41%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04074"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.45%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005952"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.35%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.368"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.70%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.391"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.12%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.02%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1398"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.87%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2500"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.84%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001272"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.05%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004373"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.27%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001199"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.01%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003741"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.40%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02601"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.09%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05816"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.49%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007558"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.69%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.57%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007306"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.70%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002013"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.84%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008825"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.10%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006979"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.87%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.27%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005790"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.36%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003527"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "55.08%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.43%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.27%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.27%"
